$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D from Excel's automatic number/type inference so that
# numeric-looking price strings (e.g. "0.780", "1.40") are kept as exact text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.659.53'
$ws.Range("E2").Value = '  +2.29%  '
$ws.Range("D3").Value = '2.076.43'
$ws.Range("E3").Value = '  +4.55%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").Value = '236.83'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("E6").Value = '  +2.48%  '
$ws.Range("D7").Value = '58.27'
$ws.Range("E7").Value = '  +8.02%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +3.75%  '
$ws.Range("D10").Value = '58.15'
$ws.Range("E10").Value = '  +1.59%  '
$ws.Range("E11").Value = '  +2.19%  '
$ws.Range("E12").Value = '  +4.09%  '
$ws.Range("D13").Value = '2.385.79'
$ws.Range("E13").Value = '  +4.35%  '
$ws.Range("D14").Value = '14.43'
$ws.Range("E14").Value = '  +2.95%  '
$ws.Range("D15").Value = '21.02'
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = '0.780'
$ws.Range("E16").Value = '  +3.85%  '
$ws.Range("D17").Value = '5.22'
$ws.Range("E17").Value = '  +3.95%  '
$ws.Range("D18").Value = '2.044.53'
$ws.Range("E18").Value = '  +1.98%  '
$ws.Range("D19").Value = '37.798.18'
$ws.Range("E20").Value = '  +21.74%  '
$ws.Range("D21").Value = '68.65'
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").Value = '224.85'
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.42'
$ws.Range("E25").Value = '  +3.05%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = '2.45'
$ws.Range("E26").Value = '  +2.23%  '
$ws.Range("D27").Value = '163.47'
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("D28").Value = '8.85'
$ws.Range("E28").Value = '  +2.91%  '
$ws.Range("E29").Value = '  +4.72%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = '1.40'
$ws.Range("E30").Value = '  +7.03%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '19.37'
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("D32").Value = '0.119'
$ws.Range("E32").Value = '  +1.49%  '
$ws.Range("E33").Value = '  +4.26%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '2.62'
$ws.Range("E34").Value = '  +14.69%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '4.48'
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("E36").Value = '  +6.43%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '3.35'
$ws.Range("E38").Value = '  +3.70%  '
$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").Value = '5.93'
$ws.Range("E39").Value = '  +12.45%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = '2.97'
$ws.Range("E41").Value = '  -2.39%  '
$ws.Range("D42").Value = '0.0961'
$ws.Range("E42").Value = '  +8.18%  '
$ws.Range("D43").Value = '1.484.34'
$ws.Range("E43").Value = '  +4.17%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0213'
$ws.Range("E44").Value = '  +4.93%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '95.34'
$ws.Range("E45").Value = '  +8.82%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Value = '4.30'
$ws.Range("E46").Value = '  +26.43%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '16.49'
$ws.Range("E47").Value = '  +11.03%  '
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("E49").Value = '  +9.56%  '
$ws.Range("E50").Value = '  +2.81%  '

# Restore column D to the workbook's default (unstyled) appearance now that
# the text values are safely stored.
$ws.Range("D2:D51").Style = "Normal"
